$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): update existing header cells A1:K1 (style s="1" already present) ---
$ws.Range("A1").Value = "productIds"
$ws.Range("B1").Value = "MSE_no_transfer"
$ws.Range("C1").Value = "MSE_transfer_basic"
$ws.Range("D1").Value = "MSE_transfer_coral"
$ws.Range("E1").Value = "MSE_transfer_sa"
$ws.Range("F1").Value = "MSE_transfer_bw"
$ws.Range("G1").Value = "MSE_transfer_nnw"
$ws.Range("H1").Value = "MAE_no_transfer"
$ws.Range("I1").Value = "MAE_transfer_basic"
$ws.Range("J1").Value = "MAE_transfer_coral"
$ws.Range("K1").Value = "MAE_transfer_sa"

# --- Header row (row 1): new header cells L1:W1, set value then copy header style from A1 ---
$ws.Range("L1").Value = "MAE_transfer_bw"
$ws.Range("M1").Value = "MAE_transfer_nnw"
$ws.Range("N1").Value = "MSE_diff_basic"
$ws.Range("O1").Value = "MSE_transfer_coral"
$ws.Range("P1").Value = "MSE_diff_sa"
$ws.Range("Q1").Value = "MSE_diff_bw"
$ws.Range("R1").Value = "MSE_diff_nnw"
$ws.Range("S1").Value = "MAE_diff_basic"
$ws.Range("T1").Value = "MAE_transfer_coral"
$ws.Range("U1").Value = "MAE_diff_sa"
$ws.Range("V1").Value = "MAE_diff_bw"
$ws.Range("W1").Value = "MAE_diff_nnw"

# Apply the bold/border/centered header style (same as A1:K1) to the newly added header cells L1:W1
$ws.Range("A1").Copy()
$ws.Range("L1:W1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-6) ---
# Row 2
$ws.Range("B2").Value = 1.301926923014102
$ws.Range("C2").Value = 0.9989431818564546
$ws.Range("D2").Value = 1.137693001037905
$ws.Range("E2").Value = 5.820702221270763
$ws.Range("F2").Value = 5.856413974995616
$ws.Range("G2").Value = 6.142760710462269
$ws.Range("H2").Value = 0.6943027395987617
$ws.Range("I2").Value = 0.6901813189040737
$ws.Range("J2").Value = 0.5502789966087731
$ws.Range("K2").Value = 0.822392764717888
$ws.Range("L2").Value = 0.808529699177241
$ws.Range("M2").Value = 0.7930880846202476
$ws.Range("N2").Value = -0.3029837411576469
$ws.Range("O2").Value = -0.1642339219761964
$ws.Range("P2").Value = 4.518775298256662
$ws.Range("Q2").Value = 4.554487051981514
$ws.Range("R2").Value = 4.840833787448168
$ws.Range("S2").Value = -0.004121420694687905
$ws.Range("T2").Value = -0.1440237429899885
$ws.Range("U2").Value = 0.1280900251191264
$ws.Range("V2").Value = 0.1142269595784794
$ws.Range("W2").Value = 0.09878534502148595

# Row 3
$ws.Range("B3").Value = 0.7858089996152173
$ws.Range("C3").Value = 1.012164272873996
$ws.Range("D3").Value = 2.777605296423082
$ws.Range("E3").Value = 5.820702221270763
$ws.Range("F3").Value = 5.856413974995616
$ws.Range("G3").Value = 6.142760710462269
$ws.Range("H3").Value = 0.5457719070059673
$ws.Range("I3").Value = 0.6718072620290712
$ws.Range("J3").Value = 0.9620219810251104
$ws.Range("K3").Value = 0.822392764717888
$ws.Range("L3").Value = 0.808529699177241
$ws.Range("M3").Value = 0.7930880846202476
$ws.Range("N3").Value = 0.226355273258779
$ws.Range("O3").Value = 1.991796296807865
$ws.Range("P3").Value = 5.034893221655546
$ws.Range("Q3").Value = 5.070604975380398
$ws.Range("R3").Value = 5.356951710847052
$ws.Range("S3").Value = 0.1260353550231038
$ws.Range("T3").Value = 0.416250074019143
$ws.Range("U3").Value = 0.2766208577119207
$ws.Range("V3").Value = 0.2627577921712737
$ws.Range("W3").Value = 0.2473161776142803

# Row 4
$ws.Range("B4").Value = 1.996350390516476
$ws.Range("C4").Value = 0.9855601843336593
$ws.Range("D4").Value = 6.135886311364529
$ws.Range("E4").Value = 5.820702221270763
$ws.Range("F4").Value = 5.856413974995616
$ws.Range("G4").Value = 6.142760710462269
$ws.Range("H4").Value = 0.8411752455834071
$ws.Range("I4").Value = 0.6803765552195126
$ws.Range("J4").Value = 1.401427603590304
$ws.Range("K4").Value = 0.822392764717888
$ws.Range("L4").Value = 0.808529699177241
$ws.Range("M4").Value = 0.7930880846202476
$ws.Range("N4").Value = -1.010790206182817
$ws.Range("O4").Value = 4.139535920848052
$ws.Range("P4").Value = 3.824351830754287
$ws.Range("Q4").Value = 3.86006358447914
$ws.Range("R4").Value = 4.146410319945794
$ws.Range("S4").Value = -0.1607986903638945
$ws.Range("T4").Value = 0.5602523580068973
$ws.Range("U4").Value = -0.01878248086551904
$ws.Range("V4").Value = -0.03264554640616601
$ws.Range("W4").Value = -0.04808716096315946

# Row 5
$ws.Range("B5").Value = 0.9698881290922886
$ws.Range("C5").Value = 0.9927559027195498
$ws.Range("D5").Value = 13.05127556167862
$ws.Range("E5").Value = 5.820702221270763
$ws.Range("F5").Value = 5.856413974995616
$ws.Range("G5").Value = 6.142760710462269
$ws.Range("H5").Value = 0.3799394289281765
$ws.Range("I5").Value = 0.4098606503903026
$ws.Range("J5").Value = 1.254331995984313
$ws.Range("K5").Value = 0.822392764717888
$ws.Range("L5").Value = 0.808529699177241
$ws.Range("M5").Value = 0.7930880846202476
$ws.Range("N5").Value = 0.02286777362726122
$ws.Range("O5").Value = 12.08138743258633
$ws.Range("P5").Value = 4.850814092178474
$ws.Range("Q5").Value = 4.886525845903327
$ws.Range("R5").Value = 5.172872581369981
$ws.Range("S5").Value = 0.02992122146212611
$ws.Range("T5").Value = 0.8743925670561368
$ws.Range("U5").Value = 0.4424533357897115
$ws.Range("V5").Value = 0.4285902702490645
$ws.Range("W5").Value = 0.4131486556920711

# Row 6
$ws.Range("B6").Value = 1.51415165254623
$ws.Range("C6").Value = 1.003750476033903
$ws.Range("D6").Value = 8.282524188105244
$ws.Range("E6").Value = 5.820702221270763
$ws.Range("F6").Value = 5.856413974995616
$ws.Range("G6").Value = 6.142760710462269
$ws.Range("H6").Value = 0.5200546449324573
$ws.Range("I6").Value = 0.4667778334864473
$ws.Range("J6").Value = 1.155958107027304
$ws.Range("K6").Value = 0.822392764717888
$ws.Range("L6").Value = 0.808529699177241
$ws.Range("M6").Value = 0.7930880846202476
$ws.Range("N6").Value = -0.5104011765123271
$ws.Range("O6").Value = 6.768372535559013
$ws.Range("P6").Value = 4.306550568724533
$ws.Range("Q6").Value = 4.342262322449385
$ws.Range("R6").Value = 4.628609057916039
$ws.Range("S6").Value = -0.05327681144601004
$ws.Range("T6").Value = 0.6359034620948468
$ws.Range("U6").Value = 0.3023381197854307
$ws.Range("V6").Value = 0.2884750542447837
$ws.Range("W6").Value = 0.2730334396877903
